# Duplicate the 15 URL rows on the "Side nav" sheet into rows 16-30,
# re-using the same shared-string values (A16=A1, A17=A2, ... A30=A15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Side nav")

for ($i = 1; $i -le 15; $i++) {
    $srcCell = $ws.Cells.Item($i, 1)
    $dstCell = $ws.Cells.Item($i + 15, 1)
    $dstCell.Value2 = $srcCell.Value2
}
